$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 438 (this shifts old rows 438..455
# down to 442..459, and keeps formatting/number-format from the row above).
$ws.Rows.Item(438).Insert()
$ws.Rows.Item(438).Insert()
$ws.Rows.Item(438).Insert()
$ws.Rows.Item(438).Insert()

# Common columns shared by every Kiwi / Vega Central Mapocho row in this block.
$commonA = 9
$commonB = "Vega Central Mapocho de Santiago"
$commonC = "Metropolitana"
$commonE = 13
$commonF = "Fruta"
$commonG = 100101
$commonH = "Berries"
$commonI = 100101007
$commonJ = "Kiwi"
$commonK = "Hayward"

$newRows = @(
    @{ Row = 438; D = 44516; L = "Especial";                M = 380; N = 11000; O = 11000; P = 11000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Curicó"; S = 1100; T = 10 },
    @{ Row = 439; D = 44516; L = "Extra (doble especial)";  M = 310; N = 12000; O = 12000; P = 12000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Curicó"; S = 1200; T = 10 },
    @{ Row = 440; D = 44516; L = "Primera";                  M = 350; N = 9000;  O = 9000;  P = 9000;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Curicó"; S = 900;  T = 10 },
    @{ Row = 441; D = 44516; L = "Segunda";                  M = 300; N = 7000;  O = 7000;  P = 7000;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Curicó"; S = 700;  T = 10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $commonA
    $ws.Cells.Item($row, 2).Value = $commonB
    $ws.Cells.Item($row, 3).Value = $commonC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $commonE
    $ws.Cells.Item($row, 6).Value = $commonF
    $ws.Cells.Item($row, 7).Value = $commonG
    $ws.Cells.Item($row, 8).Value = $commonH
    $ws.Cells.Item($row, 9).Value = $commonI
    $ws.Cells.Item($row, 10).Value = $commonJ
    $ws.Cells.Item($row, 11).Value = $commonK
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
